# #860 updated how-to-load tutorial
# Bold the "agentName" token inside the sentence
# "The agentName is a general text string giving the agent a name."
# This splits the single run into three runs: "The ", "agentName" (bold),
# " is a general text string giving the agent a name."

$d = $word.ActiveDocument

$target = "The agentName is a general text string giving the agent a name."

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq $target) {
        $para = $p
        break
    }
}

$start = $para.Range.Start
$end = $para.Range.End

# Restrict the Find to this paragraph only, so the (unrelated) earlier
# "agentName" occurrence in the document is left untouched.
$boldRange = $d.Range($start, $end)
$boldRange.Find.Execute("agentName", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$boldRange.Bold = 1
